$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 64,7
$arr[0,0] = 'DecisionTreeClassifier'
$arr[0,1] = 'CountVectorizer'
$arr[0,2] = 'Original'
$arr[0,3] = 0.8790450678460586
$arr[0,4] = 1.194809832175573
$arr[0,5] = '{"model__criterion": "entropy", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[0,6] = 29.43540979200043
$arr[1,0] = 'DecisionTreeClassifier'
$arr[1,1] = 'CountVectorizer'
$arr[1,2] = 'Stemmed'
$arr[1,3] = 0.8850742805651962
$arr[1,4] = 1.063429055611292
$arr[1,5] = '{"model__criterion": "entropy", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[1,6] = 25.61203320801724
$arr[2,0] = 'DecisionTreeClassifier'
$arr[2,1] = 'CountVectorizer'
$arr[2,2] = 'Lemmatized'
$arr[2,3] = 0.8853213582527837
$arr[2,4] = 1.042586943507194
$arr[2,5] = '{"model__criterion": "entropy", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[2,6] = 25.30355300000519
$arr[3,0] = 'DecisionTreeClassifier'
$arr[3,1] = 'CountVectorizer'
$arr[3,2] = 'Stemmed and Lemmatized'
$arr[3,3] = 0.8843010577267678
$arr[3,4] = 1.053470232089361
$arr[3,5] = '{"model__criterion": "entropy", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[3,6] = 25.30541295898729
$arr[4,0] = 'DecisionTreeClassifier'
$arr[4,1] = 'TfidfVectorizer'
$arr[4,2] = 'Original'
$arr[4,3] = 0.8803191231606777
$arr[4,4] = 2.2069451212883
$arr[4,5] = '{"model__criterion": "gini", "model__max_depth": 50, "vectorizer__max_features": null}'
$arr[4,6] = 52.07366291698418
$arr[5,0] = 'DecisionTreeClassifier'
$arr[5,1] = 'TfidfVectorizer'
$arr[5,2] = 'Stemmed'
$arr[5,3] = 0.8773545921411838
$arr[5,4] = 1.921985923250516
$arr[5,5] = '{"model__criterion": "gini", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[5,6] = 43.79587774997344
$arr[6,0] = 'DecisionTreeClassifier'
$arr[6,1] = 'TfidfVectorizer'
$arr[6,2] = 'Lemmatized'
$arr[6,3] = 0.881237611588592
$arr[6,4] = 1.896561488509178
$arr[6,5] = '{"model__criterion": "entropy", "model__max_depth": 25, "vectorizer__max_features": null}'
$arr[6,6] = 43.73308395899949
$arr[7,0] = 'DecisionTreeClassifier'
$arr[7,1] = 'TfidfVectorizer'
$arr[7,2] = 'Stemmed and Lemmatized'
$arr[7,3] = 0.8793293084365271
$arr[7,4] = 1.913800913095474
$arr[7,5] = '{"model__criterion": "gini", "model__max_depth": 50, "vectorizer__max_features": null}'
$arr[7,6] = 44.66518791698036
$arr[8,0] = 'KNeighborsClassifier'
$arr[8,1] = 'CountVectorizer'
$arr[8,2] = 'Original'
$arr[8,3] = 0.9315343353566843
$arr[8,4] = 0.3222541785240173
$arr[8,5] = '{"model__n_neighbors": 100, "model__weights": "distance", "vectorizer__max_features": null}'
$arr[8,6] = 108.4233240419999
$arr[9,0] = 'KNeighborsClassifier'
$arr[9,1] = 'CountVectorizer'
$arr[9,2] = 'Stemmed'
$arr[9,3] = 0.9427981581985472
$arr[9,4] = 0.2639323234558105
$arr[9,5] = '{"model__n_neighbors": 100, "model__weights": "distance", "vectorizer__max_features": 500}'
$arr[9,6] = 108.7184105830092
$arr[10,0] = 'KNeighborsClassifier'
$arr[10,1] = 'CountVectorizer'
$arr[10,2] = 'Lemmatized'
$arr[10,3] = 0.9384841483577351
$arr[10,4] = 0.2704619228839874
$arr[10,5] = '{"model__n_neighbors": 100, "model__weights": "distance", "vectorizer__max_features": 500}'
$arr[10,6] = 98.60888470901409
$arr[11,0] = 'KNeighborsClassifier'
$arr[11,1] = 'CountVectorizer'
$arr[11,2] = 'Stemmed and Lemmatized'
$arr[11,3] = 0.942642805890441
$arr[11,4] = 0.264656150341034
$arr[11,5] = '{"model__n_neighbors": 100, "model__weights": "distance", "vectorizer__max_features": 500}'
$arr[11,6] = 106.8641216249962
$arr[12,0] = 'KNeighborsClassifier'
$arr[12,1] = 'TfidfVectorizer'
$arr[12,2] = 'Original'
$arr[12,3] = 0.9734723514605288
$arr[12,4] = 0.3164676284790039
$arr[12,5] = '{"model__n_neighbors": 500, "model__weights": "distance", "vectorizer__max_features": null}'
$arr[12,6] = 115.3960450410086
$arr[13,0] = 'KNeighborsClassifier'
$arr[13,1] = 'TfidfVectorizer'
$arr[13,2] = 'Stemmed'
$arr[13,3] = 0.9717685256335735
$arr[13,4] = 0.2610507738590241
$arr[13,5] = '{"model__n_neighbors": 500, "model__weights": "distance", "vectorizer__max_features": null}'
$arr[13,6] = 152.728228833992
$arr[14,0] = 'KNeighborsClassifier'
$arr[14,1] = 'TfidfVectorizer'
$arr[14,2] = 'Lemmatized'
$arr[14,3] = 0.9722435836478109
$arr[14,4] = 0.2696874761581421
$arr[14,5] = '{"model__n_neighbors": 500, "model__weights": "distance", "vectorizer__max_features": null}'
$arr[14,6] = 99.70551920798607
$arr[15,0] = 'KNeighborsClassifier'
$arr[15,1] = 'TfidfVectorizer'
$arr[15,2] = 'Stemmed and Lemmatized'
$arr[15,3] = 0.9717381267377752
$arr[15,4] = 0.2474466335773468
$arr[15,5] = '{"model__n_neighbors": 500, "model__weights": "distance", "vectorizer__max_features": null}'
$arr[15,6] = 180.872806667001
$arr[16,0] = 'MultinomialNB'
$arr[16,1] = 'CountVectorizer'
$arr[16,2] = 'Original'
$arr[16,3] = 0.9531746718651161
$arr[16,4] = 0.3027505179246266
$arr[16,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[16,6] = 5.602990082988981
$arr[17,0] = 'MultinomialNB'
$arr[17,1] = 'CountVectorizer'
$arr[17,2] = 'Stemmed'
$arr[17,3] = 0.9543602927827028
$arr[17,4] = 0.2463798642158508
$arr[17,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[17,6] = 4.624882959004026
$arr[18,0] = 'MultinomialNB'
$arr[18,1] = 'CountVectorizer'
$arr[18,2] = 'Lemmatized'
$arr[18,3] = 0.95530242414735
$arr[18,4] = 0.263014413913091
$arr[18,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[18,6] = 4.775328999996418
$arr[19,0] = 'MultinomialNB'
$arr[19,1] = 'CountVectorizer'
$arr[19,2] = 'Stemmed and Lemmatized'
$arr[19,3] = 0.9544235529143363
$arr[19,4] = 0.2541276852289835
$arr[19,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[19,6] = 4.769068292021984
$arr[20,0] = 'MultinomialNB'
$arr[20,1] = 'TfidfVectorizer'
$arr[20,2] = 'Original'
$arr[20,3] = 0.964878063786201
$arr[20,4] = 0.3220246076583862
$arr[20,5] = '{"model__alpha": 0.1, "vectorizer__max_features": null}'
$arr[20,6] = 5.823270374996355
$arr[21,0] = 'MultinomialNB'
$arr[21,1] = 'TfidfVectorizer'
$arr[21,2] = 'Stemmed'
$arr[21,3] = 0.9658474053973609
$arr[21,4] = 0.2522116561730703
$arr[21,5] = '{"model__alpha": 1, "vectorizer__max_features": 500}'
$arr[21,6] = 4.76124137500301
$arr[22,0] = 'MultinomialNB'
$arr[22,1] = 'TfidfVectorizer'
$arr[22,2] = 'Lemmatized'
$arr[22,3] = 0.9659135341704523
$arr[22,4] = 0.2672414561112722
$arr[22,5] = '{"model__alpha": 1, "vectorizer__max_features": 500}'
$arr[22,6] = 4.812960500014015
$arr[23,0] = 'MultinomialNB'
$arr[23,1] = 'TfidfVectorizer'
$arr[23,2] = 'Stemmed and Lemmatized'
$arr[23,3] = 0.9658473435737553
$arr[23,4] = 0.2468320687611898
$arr[23,5] = '{"model__alpha": 1, "vectorizer__max_features": 500}'
$arr[23,6] = 4.627067333000014
$arr[24,0] = 'ComplementNB'
$arr[24,1] = 'CountVectorizer'
$arr[24,2] = 'Original'
$arr[24,3] = 0.9443289950350515
$arr[24,4] = 0.2972681760787964
$arr[24,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[24,6] = 5.396360624989029
$arr[25,0] = 'ComplementNB'
$arr[25,1] = 'CountVectorizer'
$arr[25,2] = 'Stemmed'
$arr[25,3] = 0.9438206588065693
$arr[25,4] = 0.2405405720074971
$arr[25,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[25,6] = 4.727558332990156
$arr[26,0] = 'ComplementNB'
$arr[26,1] = 'CountVectorizer'
$arr[26,2] = 'Lemmatized'
$arr[26,3] = 0.9450599365446275
$arr[26,4] = 0.2550284465154012
$arr[26,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[26,6] = 4.606934583018301
$arr[27,0] = 'ComplementNB'
$arr[27,1] = 'CountVectorizer'
$arr[27,2] = 'Stemmed and Lemmatized'
$arr[27,3] = 0.943771473528919
$arr[27,4] = 0.2375452121098836
$arr[27,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[27,6] = 4.499347666016547
$arr[28,0] = 'ComplementNB'
$arr[28,1] = 'TfidfVectorizer'
$arr[28,2] = 'Original'
$arr[28,3] = 0.9710043812036906
$arr[28,4] = 0.2907109657923381
$arr[28,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[28,6] = 5.487691458984045
$arr[29,0] = 'ComplementNB'
$arr[29,1] = 'TfidfVectorizer'
$arr[29,2] = 'Stemmed'
$arr[29,3] = 0.9695185773428346
$arr[29,4] = 0.2408114314079285
$arr[29,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[29,6] = 4.567891834012698
$arr[30,0] = 'ComplementNB'
$arr[30,1] = 'TfidfVectorizer'
$arr[30,2] = 'Lemmatized'
$arr[30,3] = 0.9698537754072606
$arr[30,4] = 0.2529229283332825
$arr[30,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[30,6] = 4.638796374987578
$arr[31,0] = 'ComplementNB'
$arr[31,1] = 'TfidfVectorizer'
$arr[31,2] = 'Stemmed and Lemmatized'
$arr[31,3] = 0.9695783102497578
$arr[31,4] = 0.2392579515775045
$arr[31,5] = '{"model__alpha": 1, "vectorizer__max_features": null}'
$arr[31,6] = 4.524278208002215
$arr[32,0] = 'LogisticRegression'
$arr[32,1] = 'CountVectorizer'
$arr[32,2] = 'Original'
$arr[32,3] = 0.9586923088015173
$arr[32,4] = 90.3775022061666
$arr[32,5] = '{"model__C": 100, "model__l1_ratio": 1.0, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[32,6] = 6068.104392333014
$arr[33,0] = 'LogisticRegression'
$arr[33,1] = 'CountVectorizer'
$arr[33,2] = 'Stemmed'
$arr[33,3] = 0.9603691484500572
$arr[33,4] = 57.12458638548851
$arr[33,5] = '{"model__C": 100, "model__l1_ratio": 0.5, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[33,6] = 3712.599043833005
$arr[34,0] = 'LogisticRegression'
$arr[34,1] = 'CountVectorizer'
$arr[34,2] = 'Lemmatized'
$arr[34,3] = 0.960091783075103
$arr[34,4] = 131.6969273602963
$arr[34,5] = '{"model__C": 100, "model__l1_ratio": 0.25, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[34,6] = 8059.46575695902
$arr[35,0] = 'LogisticRegression'
$arr[35,1] = 'CountVectorizer'
$arr[35,2] = 'Stemmed and Lemmatized'
$arr[35,3] = 0.9603303572661535
$arr[35,4] = 102.6273009749254
$arr[35,5] = '{"model__C": 10, "model__l1_ratio": 0.0, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[35,6] = 5399.201219375012
$arr[36,0] = 'LogisticRegression'
$arr[36,1] = 'TfidfVectorizer'
$arr[36,2] = 'Original'
$arr[36,3] = 0.98661685115662
$arr[36,4] = 43.29253586689632
$arr[36,5] = '{"model__C": 1, "model__l1_ratio": 0.5, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[36,6] = 2345.625237542001
$arr[37,0] = 'LogisticRegression'
$arr[37,1] = 'TfidfVectorizer'
$arr[37,2] = 'Stemmed'
$arr[37,3] = 0.9840139522250123
$arr[37,4] = 41.10230638742447
$arr[37,5] = '{"model__C": 1, "model__l1_ratio": 0.25, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[37,6] = 2198.86051054101
$arr[38,0] = 'LogisticRegression'
$arr[38,1] = 'TfidfVectorizer'
$arr[38,2] = 'Lemmatized'
$arr[38,3] = 0.9849189233537254
$arr[38,4] = 40.23901522358258
$arr[38,5] = '{"model__C": 1, "model__l1_ratio": 0.25, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[38,6] = 2193.55761495902
$arr[39,0] = 'LogisticRegression'
$arr[39,1] = 'TfidfVectorizer'
$arr[39,2] = 'Stemmed and Lemmatized'
$arr[39,3] = 0.9840052129425053
$arr[39,4] = 41.16085040569305
$arr[39,5] = '{"model__C": 1, "model__l1_ratio": 0.25, "model__penalty": "elasticnet", "model__solver": "saga", "vectorizer__max_features": null}'
$arr[39,6] = 2205.567522249999
$arr[40,0] = 'XGBClassifier'
$arr[40,1] = 'CountVectorizer'
$arr[40,2] = 'Original'
$arr[40,3] = 0.9859588231696181
$arr[40,4] = 11.99565851887067
$arr[40,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[40,6] = 831.3635886249831
$arr[41,0] = 'XGBClassifier'
$arr[41,1] = 'CountVectorizer'
$arr[41,2] = 'Stemmed'
$arr[41,3] = 0.9852839919522474
$arr[41,4] = 9.742909802993138
$arr[41,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[41,6] = 691.1724815829948
$arr[42,0] = 'XGBClassifier'
$arr[42,1] = 'CountVectorizer'
$arr[42,2] = 'Lemmatized'
$arr[42,3] = 0.9855884858753617
$arr[42,4] = 13.37820769713985
$arr[42,5] = '{"model__colsample_bytree": 0.7, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[42,6] = 918.1099537080154
$arr[43,0] = 'XGBClassifier'
$arr[43,1] = 'CountVectorizer'
$arr[43,2] = 'Stemmed and Lemmatized'
$arr[43,3] = 0.9852970436879072
$arr[43,4] = 15.2180129961835
$arr[43,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[43,6] = 1037.484150625009
$arr[44,0] = 'XGBClassifier'
$arr[44,1] = 'TfidfVectorizer'
$arr[44,2] = 'Original'
$arr[44,3] = 0.9869954725563985
$arr[44,4] = 103.6930783228742
$arr[44,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[44,6] = 6421.626396209002
$arr[45,0] = 'XGBClassifier'
$arr[45,1] = 'TfidfVectorizer'
$arr[45,2] = 'Stemmed'
$arr[45,3] = 0.9858477943893839
$arr[45,4] = 89.54037327832644
$arr[45,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[45,6] = 5565.953025166993
$arr[46,0] = 'XGBClassifier'
$arr[46,1] = 'TfidfVectorizer'
$arr[46,2] = 'Lemmatized'
$arr[46,3] = 0.9861815384593564
$arr[46,4] = 92.15152019891474
$arr[46,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 0.7, "vectorizer__max_features": null}'
$arr[46,6] = 5716.530819666979
$arr[47,0] = 'XGBClassifier'
$arr[47,1] = 'TfidfVectorizer'
$arr[47,2] = 'Stemmed and Lemmatized'
$arr[47,3] = 0.985966737897561
$arr[47,4] = 90.5555822753244
$arr[47,5] = '{"model__colsample_bytree": 0.5, "model__max_depth": 10, "model__subsample": 1, "vectorizer__max_features": null}'
$arr[47,6] = 5631.057041916996
$arr[48,0] = 'BaggingClassifier'
$arr[48,1] = 'CountVectorizer'
$arr[48,2] = 'Original'
$arr[48,3] = 0.9651625612612003
$arr[48,4] = 174.8117938488722
$arr[48,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[48,6] = 1432.970788000006
$arr[49,0] = 'BaggingClassifier'
$arr[49,1] = 'CountVectorizer'
$arr[49,2] = 'Stemmed'
$arr[49,3] = 0.9664369568545241
$arr[49,4] = 148.7001423150301
$arr[49,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[49,6] = 1223.369071166991
$arr[50,0] = 'BaggingClassifier'
$arr[50,1] = 'CountVectorizer'
$arr[50,2] = 'Lemmatized'
$arr[50,3] = 0.9644894608897217
$arr[50,4] = 134.6676911622286
$arr[50,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[50,6] = 1120.67374758402
$arr[51,0] = 'BaggingClassifier'
$arr[51,1] = 'CountVectorizer'
$arr[51,2] = 'Stemmed and Lemmatized'
$arr[51,3] = 0.9643416456508651
$arr[51,4] = 133.2518619894981
$arr[51,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[51,6] = 1115.898864250019
$arr[52,0] = 'BaggingClassifier'
$arr[52,1] = 'TfidfVectorizer'
$arr[52,2] = 'Original'
$arr[52,3] = 0.9700922613181199
$arr[52,4] = 258.0679601281881
$arr[52,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[52,6] = 2199.784869374998
$arr[53,0] = 'BaggingClassifier'
$arr[53,1] = 'TfidfVectorizer'
$arr[53,2] = 'Stemmed'
$arr[53,3] = 0.9703556166200668
$arr[53,4] = 215.0148673534393
$arr[53,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[53,6] = 1835.195166624995
$arr[54,0] = 'BaggingClassifier'
$arr[54,1] = 'TfidfVectorizer'
$arr[54,2] = 'Lemmatized'
$arr[54,3] = 0.9712287824715189
$arr[54,4] = 217.4656108409166
$arr[54,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[54,6] = 1849.731134375004
$arr[55,0] = 'BaggingClassifier'
$arr[55,1] = 'TfidfVectorizer'
$arr[55,2] = 'Stemmed and Lemmatized'
$arr[55,3] = 0.9712702706756078
$arr[55,4] = 218.3413614034653
$arr[55,5] = '{"model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[55,6] = 1857.983143707999
$arr[56,0] = 'RandomForestClassifier'
$arr[56,1] = 'CountVectorizer'
$arr[56,2] = 'Original'
$arr[56,3] = 0.9755662679587731
$arr[56,4] = 7.26015535692374
$arr[56,5] = '{"model__class_weight": "balanced", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[56,6] = 714.9583407920145
$arr[57,0] = 'RandomForestClassifier'
$arr[57,1] = 'CountVectorizer'
$arr[57,2] = 'Stemmed'
$arr[57,3] = 0.9768099552844349
$arr[57,4] = 7.271614888807139
$arr[57,5] = '{"model__class_weight": "balanced_subsample", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[57,6] = 705.0865867089597
$arr[58,0] = 'RandomForestClassifier'
$arr[58,1] = 'CountVectorizer'
$arr[58,2] = 'Lemmatized'
$arr[58,3] = 0.977316111633714
$arr[58,4] = 6.936464201907316
$arr[58,5] = '{"model__class_weight": "balanced", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[58,6] = 681.9510156660108
$arr[59,0] = 'RandomForestClassifier'
$arr[59,1] = 'CountVectorizer'
$arr[59,2] = 'Stemmed and Lemmatized'
$arr[59,3] = 0.976916426352403
$arr[59,4] = 7.077846206724644
$arr[59,5] = '{"model__class_weight": "balanced_subsample", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[59,6] = 691.1364472500281
$arr[60,0] = 'RandomForestClassifier'
$arr[60,1] = 'TfidfVectorizer'
$arr[60,2] = 'Original'
$arr[60,3] = 0.9778217950826862
$arr[60,4] = 8.81932769815127
$arr[60,5] = '{"model__class_weight": "balanced", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[60,6] = 838.9564697500318
$arr[61,0] = 'RandomForestClassifier'
$arr[61,1] = 'TfidfVectorizer'
$arr[61,2] = 'Stemmed'
$arr[61,3] = 0.9790889258763575
$arr[61,4] = 8.393872867524623
$arr[61,5] = '{"model__class_weight": "balanced_subsample", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[61,6] = 799.8681586250314
$arr[62,0] = 'RandomForestClassifier'
$arr[62,1] = 'TfidfVectorizer'
$arr[62,2] = 'Lemmatized'
$arr[62,3] = 0.9791279564943463
$arr[62,4] = 7.952092997978132
$arr[62,5] = '{"model__class_weight": "balanced", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[62,6] = 762.6157083750004
$arr[63,0] = 'RandomForestClassifier'
$arr[63,1] = 'TfidfVectorizer'
$arr[63,2] = 'Stemmed and Lemmatized'
$arr[63,3] = 0.9784742064192116
$arr[63,4] = 8.286464368800322
$arr[63,5] = '{"model__class_weight": "balanced_subsample", "model__max_depth": 100, "model__n_estimators": 100, "vectorizer__max_features": null}'
$arr[63,6] = 794.4221222919878

$ws.Range("A2:G65").Value = $arr
